$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7289664745330811
$ws.Range("B1").Value = 0.5737303495407104
$ws.Range("C1").Value = 4.71602725982666
$ws.Range("D1").Value = 2.82428765296936
$ws.Range("E1").Value = 1.207239985466003
